# Updated cryptos list on Mon May 13 17:29:14 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns for each crypto row with
# the latest scrape, and re-sorts two pairs of rows (30/31 and 50/51) whose
# relative ranking flipped (NEARProtocol/PEPE, Stellar/InjectiveProtocol).
#
# Cells in column D sometimes hold plain-looking numeric text (e.g. "7.32").
# Excel's COM layer auto-coerces such strings typed into .Value into real
# numbers, which would silently turn the text cell into a number cell. To
# keep them as text (matching the original inlineStr cells), we force the
# cell's NumberFormat to "@" (Text) immediately before assigning any D value
# that would otherwise parse as a number. Values that already contain a
# second "." (thousands-style, e.g. "62.790.83") can never parse as a plain
# number, so that trick is skipped for them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (price, volume) ; price $null means "unchanged" for that row.
$rowUpdates = @(
    @{ Row = 2;  D = "62.790.83";  E = "  +2.17%  " },
    @{ Row = 3;  D = "2.950.51";   E = $null },
    @{ Row = 4;  D = $null;        E = "  +0.06%  " },
    @{ Row = 5;  D = "593.69";     E = "  -0.28%  " },
    @{ Row = 6;  D = "148.00";     E = "  +1.93%  " },
    @{ Row = 7;  D = $null;        E = "  +0.03%  " },
    @{ Row = 8;  D = "0.508";      E = "  +1.10%  " },
    @{ Row = 9;  D = "2.945.76";   E = "  +0.44%  " },
    @{ Row = 10; D = "7.32";       E = "  +5.04%  " },
    @{ Row = 11; D = $null;        E = "  +7.16%  " },
    @{ Row = 12; D = "0.443";      E = "  +0.69%  " },
    @{ Row = 13; D = $null;        E = "  +6.25%  " },
    @{ Row = 14; D = "32.92";      E = "  -2.24%  " },
    @{ Row = 15; D = $null;        E = "  -0.82%  " },
    @{ Row = 16; D = "3.439.66";   E = "  +0.57%  " },
    @{ Row = 17; D = "62.743.78";  E = "  +2.17%  " },
    @{ Row = 18; D = "6.72";       E = "  -0.10%  " },
    @{ Row = 19; D = "2.956.93";   E = "  +0.74%  " },
    @{ Row = 20; D = "442.94";     E = "  +2.46%  " },
    @{ Row = 21; D = "13.49";      E = "  -0.10%  " },
    @{ Row = 22; D = "0.668";      E = "  -1.72%  " },
    @{ Row = 23; D = "7.04";       E = "  -1.21%  " },
    @{ Row = 24; D = "81.44";      E = "  -0.51%  " },
    @{ Row = 25; D = $null;        E = "  +2.86%  " },
    @{ Row = 26; D = "2.15";       E = "  -2.30%  " },
    @{ Row = 27; D = "11.71";      E = "  -0.76%  " },
    @{ Row = 28; D = $null;        E = "  -0.01%  " },
    @{ Row = 29; D = $null;        E = "  +0.79%  " },
    @{ Row = 32; D = "2.61";       E = "  -0.13%  " },
    @{ Row = 33; D = "26.53";      E = "  -0.66%  " },
    @{ Row = 35; D = $null;        E = "  -0.03%  " },
    @{ Row = 36; D = "0.992";      E = "  -1.95%  " },
    @{ Row = 37; D = $null;        E = "  +6.15%  " },
    @{ Row = 38; D = $null;        E = "  -0.63%  " },
    @{ Row = 39; D = "2.05";       E = "  +1.74%  " },
    @{ Row = 40; D = "49.63";      E = "  -0.46%  " },
    @{ Row = 41; D = "8.52";       E = "  -1.08%  " },
    @{ Row = 42; D = $null;        E = "  -5.47%  " },
    @{ Row = 43; D = $null;        E = "  -0.57%  " },
    @{ Row = 44; D = "39.91";      E = "  -6.45%  " },
    @{ Row = 45; D = "2.702.22";   E = "  -0.17%  " },
    @{ Row = 46; D = "135.21";     E = "  +0.71%  " },
    @{ Row = 47; D = "0.0339";     E = "  -2.53%  " },
    @{ Row = 48; D = "364.60";     E = "  -0.22%  " },
    @{ Row = 49; D = $null;        E = "  +0.02%  " }
)

foreach ($u in $rowUpdates) {
    if ($null -ne $u.D) {
        # Force text storage so numeric-looking strings (e.g. "7.32") aren't
        # silently coerced to real numbers by COM's .Value auto-typing.
        $ws.Cells.Item($u.Row, 4).NumberFormat = "@"
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($u.Row, 5).Value = $u.E
    }
}

# Rows 30/31 swapped ranking: PEPE <-> NEARProtocol.
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.26"
$ws.Range("E30").Value = "  +4.72%  "

$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0000104"
$ws.Range("E31").Value = "  +18.34%  "

# Rows 50/51 swapped ranking: InjectiveProtocol <-> Stellar.
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.105"
$ws.Range("E50").Value = "  -0.38%  "

$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.93"
$ws.Range("E51").Value = "  -3.62%  "
